# Update the "Förändrad" (Changed/Updated) date column (column C) for all
# data rows (2 through 554) from 2023-09-02 (serial 45171) to
# 2023-09-03 (serial 45172).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
$endRow = 554
$newValue = 45172

for ($r = $startRow; $r -le $endRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = $newValue
    }
}
